$wb = $excel.ActiveWorkbook

# --- Insert a new "2022-Q1" sheet right before the "总计" summary sheet ---
# Duplicate the "2021-Q4" sheet (same layout/styles as every other quarter
# sheet) and drop the copy in front of "总计", then rename + rewrite it.
$totalSheetName = "总计"
$prevQuarterSheet = $wb.Worksheets.Item("2021-Q4")
$totalSheetRef = $wb.Worksheets.Item($totalSheetName)

$prevQuarterSheet.Copy($totalSheetRef)

# NOTE: after Worksheet.Copy()/Worksheets.Add(), previously-held sheet
# references can end up repointed at the freshly created sheet instead of
# the original - always re-resolve sheets we still need by (stable) name
# after structural edits like this.
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"
$totalSheet = $wb.Worksheets.Item($totalSheetName)

# The duplicated sheet has 9 data rows (copied from 2021-Q4); 2022-Q1 only
# has 3 funds, so drop rows 5:9, leaving header + 3 data rows (A1:H4).
$newSheet.Range("A5:H9").Delete()

# Headers
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B (基金代码) and D:G (numeric-looking text fields) must stay text so
# that values like "007592" keep their leading zero, matching the source
# diff's inlineStr cells (e.g. "2.55", "94.58", ...). Re-applying the
# "Normal" cell style afterwards drops the now-unneeded "@" text-format
# styling this picks up along the way, so these land back on the same
# un-styled (s omitted) cells the diff expects.
$newSheet.Range("B2:B4").NumberFormat = "@"
$newSheet.Range("D2:G4").NumberFormat = "@"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "007592"
$newSheet.Range("C2").Value = "华夏价值精选混合"
$newSheet.Range("D2").Value = "2.55"
$newSheet.Range("E2").Value = "94.58"
$newSheet.Range("F2").Value = "9.67"
$newSheet.Range("G2").Value = "0.2466"
$newSheet.Range("H2").Value = 1

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "004648"
$newSheet.Range("C3").Value = "南方安睿混合"
$newSheet.Range("D3").Value = "8.64"
$newSheet.Range("E3").Value = "21.50"
$newSheet.Range("F3").Value = "0.72"
$newSheet.Range("G3").Value = "0.0622"
$newSheet.Range("H3").Value = 9

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "000066"
$newSheet.Range("C4").Value = "诺安鸿鑫混合"
$newSheet.Range("D4").Value = "0.74"
$newSheet.Range("E4").Value = "81.34"
$newSheet.Range("F4").Value = "3.46"
$newSheet.Range("G4").Value = "0.0256"
$newSheet.Range("H4").Value = 8

$newSheet.Range("B2:B4").Style = "Normal"
$newSheet.Range("D2:G4").Style = "Normal"

# --- Update the "总计" sheet: insert a new row for 2022-Q1 at the top of the
# data (row 2), shifting the existing rows down, and bump the A-column index
# counters so they stay 0-based sequential.
$totalSheet.Rows.Item(2).Insert()

# Row insertion drags inherited formatting from the row above (bold header
# style) onto B2:D2 - clear that so the new data row matches the plain style
# used by every other data row, then restore the dedicated index-column
# style (matching A3:A7) onto A2.
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.33

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
